# New crime data collected — weekly CompStat update.
# Updates the "Volume/Number" and "Report Covering the Week" header text,
# plus refreshed crime-complaint figures for rows 15-27 (Rape .. Other Sex
# Crimes) on the single worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Header rich text: bump the report "Number" and shift both dates a
#    week forward (the cells hold multi-run rich text, so only the
#    specific run's characters are replaced, leaving the rest intact).
# ---------------------------------------------------------------------

$volCell = $ws.Range("A8")
$volText = $volCell.Value2
$volIdx = $volText.LastIndexOf("6")
$volCell.Characters($volIdx + 1, 1).Text = "7"

$weekCell = $ws.Range("C9")
$weekText = $weekCell.Value2
$d1 = "2/5/2024"
$d1Idx = $weekText.IndexOf($d1)
$weekCell.Characters($d1Idx + 1, $d1.Length).Text = "2/12/2024"

$weekText2 = $weekCell.Value2
$d2 = "2/11/2024"
$d2Idx = $weekText2.IndexOf($d2)
$weekCell.Characters($d2Idx + 1, $d2.Length).Text = "2/18/2024"

# ---------------------------------------------------------------------
# 2) Plain value refreshes: cells whose style/type does not change,
#    only the underlying number.
# ---------------------------------------------------------------------

$simpleValues = @{
    "J15" = 4
    "K15" = -50
    "N15" = -71.428571428571
    "C16" = 3
    "E16" = 50
    "F16" = 11
    "G16" = 12
    "H16" = -8.333333333333
    "I16" = 19
    "J16" = 18
    "K16" = 5.555555555555
    "M16" = 58.333333333333
    "N16" = -82.727272727272
    "C17" = 1
    "E17" = -66.666666666666
    "F17" = 12
    "G17" = 15
    "H17" = -20
    "I17" = 29
    "J17" = 25
    "K17" = 16
    "L17" = 45
    "M17" = 123.076923076923
    "N17" = -60.273972602739
    "C18" = 6
    "D18" = 8
    "E18" = -25
    "F18" = 20
    "G18" = 20
    "H18" = 0
    "I18" = 26
    "J18" = 41
    "K18" = -36.585365853658
    "L18" = -38.095238095238
    "M18" = -43.478260869565
    "N18" = -80.597014925373
    "D19" = 13
    "E19" = -15.384615384615
    "F19" = 53
    "H19" = -19.696969696969
    "I19" = 101
    "J19" = 109
    "K19" = -7.339449541284
    "L19" = -14.406779661016
    "M19" = -8.181818181818
    "N19" = -45.989304812834
    "J20" = 4
    "K20" = -25
    "L20" = -57.142857142857
    "M20" = -40
    "N20" = -95.522388059701
    "C21" = 21
    "D21" = 29
    "E21" = -27.586206896551
    "F21" = 99
    "G21" = 118
    "H21" = -16.101694915254
    "I21" = 180
    "J21" = 201
    "K21" = -10.447761194029
    "L21" = -21.052631578947
    "M21" = -3.743315508021
    "N21" = -68.965517241379
    "C23" = 2
    "D23" = 2
    "E23" = 0
    "F23" = 9
    "G23" = 8
    "H23" = 12.5
    "I23" = 10
    "J23" = 15
    "K23" = -33.333333333333
    "L23" = -44.444444444444
    "M23" = -33.333333333333
    "C24" = 24
    "D24" = 30
    "E24" = -20
    "F24" = 99
    "G24" = 111
    "H24" = -10.810810810810
    "I24" = 165
    "J24" = 184
    "K24" = -10.326086956521
    "L24" = -34.523809523809
    "M24" = -14.507772020725
    "C25" = 4
    "E25" = -42.857142857142
    "F25" = 24
    "H25" = -17.241379310344
    "I25" = 51
    "J25" = 56
    "K25" = -8.928571428571
    "L25" = -13.559322033898
    "M25" = -10.526315789473
    "G26" = 6
    "H26" = -83.333333333333
    "J26" = 8
    "K26" = -62.5
    "G27" = 2
    "H27" = 250
    "I27" = 9
    "J27" = 5
    "K27" = 80
    "L27" = -40
}

foreach ($addr in $simpleValues.Keys) {
    $ws.Range($addr).Value = $simpleValues[$addr]
}

# ---------------------------------------------------------------------
# 3) Cells that flip between the "no data" text placeholder style
#    (right-aligned text showing "0" / "***.*") and the normal numeric
#    style. A few weeks ago these precincts had zero reports (hence the
#    placeholder text); this week they have real counts, and vice versa
#    for a couple of cells. Since the placeholder text re-uses shared
#    strings, we copy a still-unchanged donor cell (row 22, which keeps
#    its placeholders this week) to pick up the exact style + shared
#    string, then overwrite with the real value where needed.
# ---------------------------------------------------------------------

# Donor cells (row 22 is untouched by this week's refresh):
#   C22 -> style 14, shared string "0"      (text placeholder for 0)
#   F22 -> style 16                         (plain numeric style)
#   H22 -> style 15                         (percent-change numeric style)

# Row 15 (Rape)
$ws.Range("C22").Copy($ws.Range("C15"))      # -> text "0"
$ws.Range("F22").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1
$ws.Range("H22").Copy($ws.Range("E15"))
$ws.Range("E15").Value = -100
$ws.Range("H22").Copy($ws.Range("M15"))
$ws.Range("M15").Value = 100

# Row 20 (G.L.A.)
$ws.Range("C22").Copy($ws.Range("C20"))      # -> text "0"
$ws.Range("F22").Copy($ws.Range("D20"))
$ws.Range("D20").Value = 2
$ws.Range("H22").Copy($ws.Range("E20"))
$ws.Range("E20").Value = -100
$ws.Range("F22").Copy($ws.Range("G20"))
$ws.Range("G20").Value = 2
$ws.Range("H22").Copy($ws.Range("H20"))
$ws.Range("H20").Value = 0

# Row 26 (UCR Rape*)
$ws.Range("C22").Copy($ws.Range("C26"))      # -> text "0"
$ws.Range("F22").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 1
$ws.Range("H22").Copy($ws.Range("E26"))
$ws.Range("E26").Value = -100
